$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.123.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.37%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.601.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.55%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'603.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'196.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.26%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.29%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.11%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.207"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.51%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.75%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'53.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.44%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0000304"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.14%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'9.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.16%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.171.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.68%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'13.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.07%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'591.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.56%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'70.257.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.29%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'Chainlink"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'19.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.70%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'WrappedEther"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'3.596.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.49%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E21").Value = "'  +0.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'17.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.69%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'101.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.10%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.04%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.92%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D29").Value = "'33.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.72%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.51%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.02%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -2.96%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.23%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'63.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.15%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0₃0891"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +8.03%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.948.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.42%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.59%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'523.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.64%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.03%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'36.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.11%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.60%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.31%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.58%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0455"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.72%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.73%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.57%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.140"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.34%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.98%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.21%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.000253"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.06%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +3.35%  "
$ws.Range("E51").Style = "Normal"
